# Generate Report for Handoff
# Updates the "Status" text from "Handed back: in sync with en-US" to
# "Ready for handoff", refreshes the associated timestamps, and narrows
# the Status-related columns that no longer need to fit the long string.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# Target stored column width (OOXML) is 17.2159881591797 characters. The
# COM ColumnWidth setter here quantizes to increments of 1/6 of a
# character (i.e. stored = round(chars*6)/6 + 5/6), so 16.333333333333332
# is the input that lands on the closest achievable grid value,
# 17.166666666666668.
$narrowWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-05 19:06:41"

$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-05 19:06:36"

$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
